$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 178.1
$ws.Range("I33").Value = 195.375
$ws.Range("K33").Value = 195.375
$ws.Range("M33").Value = 33.625
$ws.Range("H113").Value = 5168.5
$ws.Range("I113").Value = 4833
$ws.Range("K113").Value = 4833
$ws.Range("M113").Value = -1579
$ws.Range("H124").Value = 149775
$ws.Range("J124").Value = 149775
$ws.Range("L124").Value = 149775
$ws.Range("N124").Value = -159595
$ws.Range("H131").Value = 564012.5
$ws.Range("I131").Value = 722765.9
$ws.Range("K131").Value = 2168297.7
$ws.Range("M131").Value = -2163257.7
$ws.Range("H132").Value = 4603.353
$ws.Range("I132").Value = 4712.5454
$ws.Range("K132").Value = 14137.6362
$ws.Range("M132").Value = -11607.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3595.8928
$ws.Range("I32").Value = 3618.037
$ws.Range("K32").Value = 3618.037
$ws.Range("M32").Value = -3331.037
$ws.Range("H61").Value = 2671
$ws.Range("I61").Value = 1257
$ws.Range("K61").Value = 1257
$ws.Range("M61").Value = -1045
$ws.Range("H74").Value = 181231.03
$ws.Range("I74").Value = 187175.4
$ws.Range("K74").Value = 187175.4
$ws.Range("M74").Value = -186301.4
$ws.Range("H77").Value = 181231.03
$ws.Range("I77").Value = 187175.4
$ws.Range("K77").Value = 935877
$ws.Range("M77").Value = -931509
$ws.Range("H86").Value = 90000
$ws.Range("J86").Value = 90000
$ws.Range("L86").Value = 90000
$ws.Range("N86").Value = -92372
$ws.Range("H89").Value = 90000
$ws.Range("J89").Value = 90000
$ws.Range("L89").Value = 270000
$ws.Range("N89").Value = -281856
$ws.Range("H102").Value = 2781.8572
$ws.Range("I102").Value = 2370.5
$ws.Range("K102").Value = 2370.5
$ws.Range("M102").Value = -748.5
$ws.Range("H136").Value = 2671
$ws.Range("I136").Value = 1257
$ws.Range("K136").Value = 3771
$ws.Range("M136").Value = -1221

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3468.4092
$ws.Range("I86").Value = 3057.375
$ws.Range("J86").Value = 4564.5
$ws.Range("K86").Value = 3057.375
$ws.Range("L86").Value = 4564.5
$ws.Range("M86").Value = -1934.375
$ws.Range("N86").Value = -6810.5
$ws.Range("H89").Value = 3468.4092
$ws.Range("I89").Value = 3057.375
$ws.Range("J89").Value = 4564.5
$ws.Range("K89").Value = 15286.875
$ws.Range("L89").Value = 22822.5
$ws.Range("M89").Value = -9670.875
$ws.Range("N89").Value = -34054.5
$ws.Range("H134").Value = 1843.9131
$ws.Range("I134").Value = 1235.0883
$ws.Range("K134").Value = 3705.2649
$ws.Range("M134").Value = -1170.2649

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4044.547
$ws.Range("I31").Value = 3033.6897
$ws.Range("K31").Value = 3033.6897
$ws.Range("M31").Value = -2738.6897
$ws.Range("H34").Value = 4044.547
$ws.Range("I34").Value = 3033.6897
$ws.Range("K34").Value = 3033.6897
$ws.Range("M34").Value = -2831.6897
$ws.Range("H43").Value = 41999.5
$ws.Range("J43").Value = 41999.5
$ws.Range("L43").Value = 41999.5
$ws.Range("N43").Value = -42367.5
$ws.Range("H101").Value = 41999.5
$ws.Range("J101").Value = 41999.5
$ws.Range("L101").Value = 41999.5
$ws.Range("N101").Value = -48489.5
$ws.Range("H122").Value = 2059.4614
$ws.Range("I122").Value = 2142.6
$ws.Range("K122").Value = 6427.799999999999
$ws.Range("M122").Value = -3977.799999999999
$ws.Range("H132").Value = 12825713
$ws.Range("J132").Value = 27783878
$ws.Range("L132").Value = 83351634
$ws.Range("N132").Value = -83356694
$ws.Range("H134").Value = 4350.55
$ws.Range("I134").Value = 4726.0625
$ws.Range("K134").Value = 14178.1875
$ws.Range("M134").Value = -11643.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 73.81
$ws.Range("I4").Value = 74.35354
$ws.Range("K4").Value = 223.06062
$ws.Range("M4").Value = -111.06062
$ws.Range("H107").Value = 1103.091
$ws.Range("J107").Value = 1293.8572
$ws.Range("L107").Value = 3881.5716
$ws.Range("N107").Value = -7721.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 5594938
$ws.Range("J15").Value = 44305.375
$ws.Range("L15").Value = 44305.375
$ws.Range("N15").Value = -44881.375
$ws.Range("H81").Value = 5594938
$ws.Range("J81").Value = 44305.375
$ws.Range("L81").Value = 44305.375
$ws.Range("N81").Value = -46301.375
$ws.Range("H84").Value = 5594938
$ws.Range("J84").Value = 44305.375
$ws.Range("L84").Value = 132916.125
$ws.Range("N84").Value = -142900.125
$ws.Range("H113").Value = 2585.8333
$ws.Range("I113").Value = 2496.3333
$ws.Range("K113").Value = 2496.3333
$ws.Range("M113").Value = -326.3332999999998
$ws.Range("H122").Value = 7697109
$ws.Range("I122").Value = 12824344
$ws.Range("K122").Value = 38473032
$ws.Range("M122").Value = -38470582
$ws.Range("H132").Value = 3333.3333
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1699.8572
$ws.Range("I10").Value = 533.3333
$ws.Range("J10").Value = 2574.75
$ws.Range("K10").Value = 533.3333
$ws.Range("L10").Value = 2574.75
$ws.Range("M10").Value = -393.3333
$ws.Range("N10").Value = -2854.75
$ws.Range("H122").Value = 13152.637
$ws.Range("J122").Value = 13835
$ws.Range("L122").Value = 41505
$ws.Range("N122").Value = -46405
$ws.Range("H131").Value = 63750
$ws.Range("I131").Value = 62857.145
$ws.Range("J131").Value = 70000
$ws.Range("K131").Value = 62857.145
$ws.Range("L131").Value = 70000
$ws.Range("M131").Value = -57817.145
$ws.Range("N131").Value = -80080
$ws.Range("H136").Value = 5067.7827
$ws.Range("I136").Value = 4007.5
$ws.Range("J136").Value = 6224.4546
$ws.Range("K136").Value = 12022.5
$ws.Range("L136").Value = 18673.3638
$ws.Range("M136").Value = -9472.5
$ws.Range("N136").Value = -23773.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 39676
$ws.Range("I70").Value = 34514.5
$ws.Range("K70").Value = 34514.5
$ws.Range("M70").Value = -34199.5
$ws.Range("H73").Value = 39676
$ws.Range("I73").Value = 34514.5
$ws.Range("K73").Value = 34514.5
$ws.Range("M73").Value = -33422.5
$ws.Range("H132").Value = 2132.8823
$ws.Range("I132").Value = 1842.138
$ws.Range("K132").Value = 5526.414
$ws.Range("M132").Value = -2996.414
